# Ajout du niveau 4 dans le fichier excel des niveaux
$wb = $excel.ActiveWorkbook

$wsNiveaux = $wb.Worksheets.Item("Niveau 1-5")
$wsTous = $wb.Worksheets.Item("Tous_les_niveaux")

# --- New "Niveau 4" header block (row 26), styled like the existing
# "Niveau 1" / "Niveau 2" / "Niveau 3" headers (row 1 / row 9 / row 17). ---
$wsNiveaux.Range("C9:G9").Copy() | Out-Null
$wsNiveaux.Range("C26:G26").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$wsNiveaux.Range("C26").Value = "Niveau 4"
$wsNiveaux.Range("C26:G26").Merge() | Out-Null

# --- New "face" glyph block (rows 27:34), identical pattern to the one
# already drawn on the "Tous_les_niveaux" summary sheet at I16:M23. ---
$wsTous.Range("I16:M23").Copy() | Out-Null
$wsNiveaux.Range("C27:G34").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$wsTous.Range("I16:M23").Copy() | Out-Null
$wsNiveaux.Range("C27:G34").PasteSpecial(-4163) | Out-Null  # xlPasteValues

$excel.CutCopyMode = 0

# --- Selections / active sheet to match the saved state. ---
$wsTous.Activate() | Out-Null
$wsTous.Range("I16:M23").Select() | Out-Null

$wsNiveaux.Activate() | Out-Null
$wsNiveaux.Range("L31").Select() | Out-Null
